$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; existing rows 19:93 shift down to 20:94.
$ws.Rows("19:19").Insert()

# Populate the newly-inserted row 19 with the new weekly price record.
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 44764
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100112038
$ws.Range("G19").Value = "Cebollín baby"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 2750
$ws.Range("N19").Value = "`$/paquete 1,5 a 2 kilos"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 1375
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = "Hortaliza"
